# The "GitHub URL" table row currently holds a hyperlink field
# ("w1-set-exercises-mahaalmas-main") followed by a separate "/" run.
# Replace that whole cell paragraph with a single plain run containing
# the full URL as literal text (no hyperlink, no leftover run formatting).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the row whose first cell reads "GitHub URL" (row 3 in the
# current layout) so the script still finds the right cell even if the
# table's row order ever changes.
$targetRow = 3
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    if ($t.Cell($i, 1).Range.Text.Trim() -eq "GitHub URL") {
        $targetRow = $i
        break
    }
}

$cell = $t.Cell($targetRow, 2)
$rng = $cell.Range

$newXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t>https://github.com/mahaalmas/w1-set-exercises-mahaalmas-main</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# InsertXML replaces the full contents of the target range (field,
# hyperlink run and the trailing "/" run included) with the plain-text
# run described above. Cast to [void] so nothing lands on the output
# stream.
[void]$rng.InsertXML($newXml)
